$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for A1:H31 after the role-change schedule fix.
# (Test1/test2 columns collapsed into a single corrected "test1(F)" column,
#  and several shift values were corrected across the sheet.)
$data = @(
  @("Date", "Helen(F)", "Lili(F)", "Matthew(F)", "Ka(F)", "Kit(F)", "Paul(F)", "test1(F)"),
  @("17/03/2025", "7-16", "15-24", "off", "off", "off", "10-19", "15-24"),
  @("18/03/2025", "10-19", "off", "15-24", "off", "15-24", "7-16", "off"),
  @("19/03/2025", "off", "15-24", "15-24", "off", "7-16", "off", "10-19"),
  @("20/03/2025", "15-24", "10-19", "off", "off", "off", "15-24", "7-16"),
  @("21/03/2025", "off", "7-16", "off", "10-19", "15-24", "15-24", "off"),
  @("22/03/2025", "off", "off", "off", "7-16", "10-19", "off", "15-24"),
  @("23/03/2025", "off", "off", "7-16", "off", "off", "10-19", "15-24"),
  @("24/03/2025", "off", "7-16", "10-19", "15-24", "15-24", "off", "off"),
  @("25/03/2025", "15-24", "10-19", "off", "15-24", "off", "off", "7-16"),
  @("26/03/2025", "15-24", "7-16", "off", "off", "15-24", "off", "10-19"),
  @("27/03/2025", "off", "15-24", "off", "10-19", "off", "7-16", "15-24"),
  @("28/03/2025", "15-24", "15-24", "off", "off", "7-16", "10-19", "off"),
  @("29/03/2025", "15-24", "10-19", "off", "off", "off", "7-16", "off"),
  @("30/03/2025", "10-19", "off", "7-16", "off", "off", "off", "off"),
  @("31/03/2025", "10-19", "15-24", "15-24", "off", "off", "7-16", "off"),
  @("01/04/2025", "7-16", "15-24", "15-24", "off", "off", "10-19", "off"),
  @("02/04/2025", "15-24", "15-24", "7-16", "off", "off", "10-19", "off"),
  @("03/04/2025", "7-16", "off", "off", "15-24", "10-19", "15-24", "off"),
  @("04/04/2025", "7-16", "off", "10-19", "off", "15-24", "15-24", "off"),
  @("05/04/2025", "off", "15-24", "off", "7-16", "10-19", "off", "off"),
  @("06/04/2025", "10-19", "15-24", "off", "7-16", "off", "off", "off"),
  @("07/04/2025", "off", "15-24", "10-19", "7-16", "off", "15-24", "off"),
  @("08/04/2025", "10-19", "15-24", "off", "7-16", "15-24", "off", "off"),
  @("09/04/2025", "15-24", "15-24", "10-19", "7-16", "off", "off", "off"),
  @("10/04/2025", "off", "15-24", "10-19", "7-16", "off", "15-24", "off"),
  @("11/04/2025", "10-19", "15-24", "15-24", "7-16", "off", "off", "off"),
  @("12/04/2025", "off", "off", "15-24", "7-16", "off", "10-19", "off"),
  @("13/04/2025", "10-19", "off", "off", "15-24", "off", "7-16", "off"),
  @("14/04/2025", "off", "15-24", "10-19", "15-24", "off", "7-16", "off"),
  @("15/04/2025", "off", "15-24", "10-19", "15-24", "off", "7-16", "off")
)

# Column A holds dd/mm/yyyy text like "01/04/2025" — force text format first so
# Excel doesn't auto-convert the ambiguous (day<=12) ones into date serials.
$ws.Range("A2:A31").NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

# Column I (the old "test2(F)" helper column) no longer exists after the fix;
# remove its leftover data so the used range shrinks back to A1:H31.
$ws.Range("I1:I31").Delete()
